$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '49.759.65'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +3.04%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.611.78'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +4.15%  '
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '109.96'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +1.37%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '323.28'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.51%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.534'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').Value = '  -0.08%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.563'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +3.68%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '41.00'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +2.67%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '20.61'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.28%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.0823'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('E13').Value = '  +0.70%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +1.95%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.018.87'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +4.05%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.599.86'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +3.73%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.870'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +3.00%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '49.741.21'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +3.35%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '3.13'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +13.34%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '13.34'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('E21').Value = '  +0.56%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.0₃0952'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.49%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '283.03'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.46%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '72.76'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('E25').Value = '  +0.24%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '26.71'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +3.42%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -6.76%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '0.145'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +3.99%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '9.95'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.41%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '36.17'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.21%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '49.56'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.92%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '19.68'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.36%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.47'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +1.98%  '
$ws.Range('E35').Value = '  -0.06%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.0793'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('E37').Value = '  +4.78%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.73'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +1.58%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.06'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +4.03%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '22.84'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +6.48%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '123.52'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +1.78%  '
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('E45').Value = '  +4.86%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.046.96'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.32%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.03'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +9.56%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.17'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +9.28%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '9.07'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +0.90%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '5.37'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.93%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '81.93'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +2.24%  '
